# Leave Card update (4/12/2023 4:43 PM)
# A new leave entry (SP(1-0-0)) is inserted as row 26 on the "Sheet1" leave
# schedule table, pushing the existing period rows down by one (through the
# row that used to be the last data row, row 132, which becomes row 133).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Duplicate the final table row (132) down to the new final row (133)
#    so the period schedule has room for the newly inserted entry.
#    Force the row into existence first (PasteSpecial needs an existing
#    cell to target), then copy every column's value+format down, and
#    finally restore the calculated "EARNED " column formula.
# ---------------------------------------------------------------------
$ws.Range("A133").Value = 0
$lastRowCols = @("A","B","C","D","E","F","H","I","J","K")
foreach ($col in $lastRowCols) {
    $ws.Range($col + "132").Copy()
    $ws.Range($col + "133").PasteSpecial($xlPasteFormats)
}
$ws.Range("G132").Copy()
$ws.Range("G133").PasteSpecial($xlPasteFormats)
$ws.Range("G133").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("A133").ClearContents()

# ---------------------------------------------------------------------
# 2) Re-format the (old) row 132 so it matches the regular/blank period
#    row style used throughout the table (it used to be styled as the
#    special closing row, which now lives at row 133).
# ---------------------------------------------------------------------
foreach ($col in $lastRowCols) {
    $ws.Range($col + "131").Copy()
    $ws.Range($col + "132").PasteSpecial($xlPasteFormats)
}
$ws.Range("G131").Copy()
$ws.Range("G132").PasteSpecial($xlPasteFormats)
$ws.Range("G132").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# ---------------------------------------------------------------------
# 3) Shift the PERIOD date column down by one row, from row 103 back up
#    to row 28, so every period date moves to the row below it.
# ---------------------------------------------------------------------
for ($r = 103; $r -ge 28; $r--) {
    $prevValue = $ws.Cells.Item($r - 1, 1).Value()
    $ws.Cells.Item($r, 1).Value = $prevValue
}

# ---------------------------------------------------------------------
# 4) Row 27 takes the date that used to be in row 26, and its
#    PARTICULARS cell format changes to match the plain data style.
# ---------------------------------------------------------------------
$ws.Range("A27").Value = 45017
$ws.Range("C27").Copy()
$ws.Range("B27").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# 5) Row 26 becomes the newly inserted leave entry: no PERIOD date,
#    PARTICULARS = "SP(1-0-0)", and a BALANCE-column date of 4/5/2023
#    (serial 44999) in column K.
# ---------------------------------------------------------------------
$ws.Range("A26").ClearContents()
$ws.Range("B26").Value = "SP(1-0-0)"
$ws.Range("K24").Copy()
$ws.Range("K26").PasteSpecial($xlPasteFormats)
$ws.Range("K26").Value = 44999

# ---------------------------------------------------------------------
# 6) Grow Table1 so the table range covers the newly added row 133.
# ---------------------------------------------------------------------
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A8:K133"))

Write-Host "Leave card row inserted and table resized."
